# Generate Report for Archive
#
# 1. Update the localization "Status" text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F2, zh-cn!C2,
#    de-de!C2 all share the same string).
# 2. Narrow the "Status" column(s) on the Overview sheet (E:F) and on the
#    per-language sheets (C) to their new, smaller width.

$wb = $excel.ActiveWorkbook

$statusText = "In Translation"
$newStatusWidth = 12.5   # renders as the narrower ~13.41-wide column

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = $newStatusWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Columns.Item(3).ColumnWidth = $newStatusWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Columns.Item(3).ColumnWidth = $newStatusWidth
